# 294-...-Newcreateloan.xlsx — "Loan RBI, Variable Instalments"
#
# The "Repayment Schedule" sheet gains a new (currently unused/blank) column
# between the existing "In Advance" (M) and "Late" (old N) columns, so the
# schedule can later carry a "Variable Instalments" style value. This pushes
# the old N ("Late") to O and the old P ("Outstanding") to Q, widening the
# sheet from A:P to A:Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N, shifting Late/Outstanding one column right.
$ws.Columns("N:N").Insert()

# New column gets an explicit (non bestFit) width.
$ws.Columns("N:N").ColumnWidth = 9.14

# Park the cursor/selection where the author last left it.
[void]$ws.Range("T5").Select()
